# Regenerate the "K" column (column G) values on the active worksheet.
# The new values replace the previous "Strike#" simulation output with a
# freshly generated set of K values (std/mean recalculated, s_vals written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    8  = 2
    9  = 2
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 2
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 2
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 0
    30 = 2
    31 = 2
    32 = 1
    33 = 0
    34 = 0
    35 = 1
    36 = 3
    37 = 0
    38 = 1
    39 = 0
    40 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
